$wb = $excel.ActiveWorkbook

# Sheet ALC, row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 864.87933
$ws.Range("I15").Value = 864.87933
$ws.Range("K15").Value = 2594.63799
$ws.Range("M15").Value = -2425.63799

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2015.3846
$ws.Range("I40").Value = 1836.6666
$ws.Range("J40").Value = 2259.0908
$ws.Range("K40").Value = 1836.6666
$ws.Range("L40").Value = 2259.0908
$ws.Range("M40").Value = -1661.6666
$ws.Range("N40").Value = -2609.0908

# Sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 240.05
$ws.Range("I53").Value = 161.61539
$ws.Range("J53").Value = 385.7143
$ws.Range("K53").Value = 161.61539
$ws.Range("L53").Value = 385.7143
$ws.Range("M53").Value = 475.38461
$ws.Range("N53").Value = -1659.7143

# Sheet ALC, row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 83337510
$ws.Range("I64").Value = 250001870
$ws.Range("J64").Value = 5328.75
$ws.Range("K64").Value = 250001870
$ws.Range("L64").Value = 5328.75
$ws.Range("M64").Value = -250001622
$ws.Range("N64").Value = -5824.75

# Sheet ALC, row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 83337510
$ws.Range("I67").Value = 250001870
$ws.Range("J67").Value = 5328.75
$ws.Range("K67").Value = 250001870
$ws.Range("L67").Value = 5328.75
$ws.Range("M67").Value = -250001012
$ws.Range("N67").Value = -7044.75

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 52022
$ws.Range("I100").Value = 72680
$ws.Range("J100").Value = 3820
$ws.Range("K100").Value = 72680
$ws.Range("L100").Value = 3820
$ws.Range("M100").Value = -72139
$ws.Range("N100").Value = -4902

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 17489.834
$ws.Range("J112").Value = 25839.75
$ws.Range("L112").Value = 77519.25
$ws.Range("N112").Value = -79735.25

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1592.1818
$ws.Range("J129").Value = 1746
$ws.Range("L129").Value = 5238
$ws.Range("N129").Value = -15238

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1814.5
$ws.Range("I45").Value = 1540.4
$ws.Range("J45").Value = 2271.3333
$ws.Range("K45").Value = 1540.4
$ws.Range("L45").Value = 2271.3333
$ws.Range("M45").Value = -1163.4
$ws.Range("N45").Value = -3025.3333

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2030
$ws.Range("I122").Value = 2237.3333
$ws.Range("J122").Value = 1905.6
$ws.Range("K122").Value = 6711.999899999999
$ws.Range("L122").Value = 5716.799999999999
$ws.Range("M122").Value = -4261.999899999999
$ws.Range("N122").Value = -10616.8

# Sheet BSM, row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 64.89655
$ws.Range("I7").Value = 31.428572
$ws.Range("J7").Value = 96.13333
$ws.Range("K7").Value = 31.428572
$ws.Range("L7").Value = 96.13333
$ws.Range("M7").Value = 81.571428
$ws.Range("N7").Value = -322.13333

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 414.55554
$ws.Range("I22").Value = 480.16666
$ws.Range("J22").Value = 283.33334
$ws.Range("K22").Value = 480.16666
$ws.Range("L22").Value = 283.33334
$ws.Range("M22").Value = -130.16666
$ws.Range("N22").Value = -983.33334

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 58832240
$ws.Range("I132").Value = 83343930
$ws.Range("J132").Value = 4191.6
$ws.Range("K132").Value = 250031790
$ws.Range("L132").Value = 12574.8
$ws.Range("M132").Value = -250029260
$ws.Range("N132").Value = -17634.8

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 719.18335
$ws.Range("J131").Value = 916.97437
$ws.Range("L131").Value = 2750.92311
$ws.Range("N131").Value = -12830.92311

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1100
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9900
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14960

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1427.4445
$ws.Range("I122").Value = 1393.375
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 4180.125
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -1730.125
$ws.Range("N122").Value = -10000

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2043.375
$ws.Range("I46").Value = 2377.8
$ws.Range("J46").Value = 1486
$ws.Range("K46").Value = 2377.8
$ws.Range("L46").Value = 1486
$ws.Range("M46").Value = -2189.8
$ws.Range("N46").Value = -1862

# Sheet LTW, row 64
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 32700
$ws.Range("J64").Value = 32700
$ws.Range("L64").Value = 32700
$ws.Range("N64").Value = -33150

# Sheet LTW, row 67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 32700
$ws.Range("J67").Value = 32700
$ws.Range("L67").Value = 32700
$ws.Range("N67").Value = -34260

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1680.3889
$ws.Range("I68").Value = 1570.138
$ws.Range("J68").Value = 2137.1428
$ws.Range("K68").Value = 1570.138
$ws.Range("L68").Value = 2137.1428
$ws.Range("M68").Value = -821.1379999999999
$ws.Range("N68").Value = -3635.1428

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1680.3889
$ws.Range("I71").Value = 1570.138
$ws.Range("J71").Value = 2137.1428
$ws.Range("K71").Value = 7850.69
$ws.Range("L71").Value = 10685.714
$ws.Range("M71").Value = -4106.69
$ws.Range("N71").Value = -18173.714

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 900.2727
$ws.Range("I82").Value = 900.2727
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 900.2727
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -539.2727
$ws.Range("N82").ClearContents()

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 900.2727
$ws.Range("I85").Value = 900.2727
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 900.2727
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 347.7273
$ws.Range("N85").ClearContents()

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 65498.25
$ws.Range("I122").Value = 93668.91
$ws.Range("J122").Value = 3522.8
$ws.Range("K122").Value = 281006.73
$ws.Range("L122").Value = 10568.4
$ws.Range("M122").Value = -278556.73
$ws.Range("N122").Value = -15468.4

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6584.4116
$ws.Range("I132").Value = 7162
$ws.Range("K132").Value = 21486
$ws.Range("M132").Value = -18956

# Sheet WVR, row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 9000
$ws.Range("J47").Value = 9000
$ws.Range("L47").Value = 9000
$ws.Range("N47").Value = -10144

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 45466540
$ws.Range("I132").Value = 55569270
$ws.Range("J132").Value = 4252
$ws.Range("K132").Value = 166707810
$ws.Range("L132").Value = 12756
$ws.Range("M132").Value = -166705280
$ws.Range("N132").Value = -17816
